# Updates the "adjudication-date" StructureDefinition workbook to the
# newer IG publication run (v6.0.0, 2022-01-21), matching the upstream
# Alvearie gh-pages deploy commit.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------

# Version bump
$metadata.Range("B3").Value = "6.0.0"

# Publication date bump
$metadata.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now populated
$metadata.Range("B9").Value = "Alvearie Team"

# Row 10 was a duplicated "Contact" row (shared with row 11); turn it
# into the real "Jurisdiction" row ...
$metadata.Range("A10").Value = "Jurisdiction"
$metadata.Range("B10").Value = "United States of America"

# ... and drop the now-redundant duplicate "Contact" row entirely,
# shifting everything below it up by one (21 rows -> 20 rows).
$metadata.Range("A11:B11").EntireRow.Delete()

# --- Elements sheet ---------------------------------------------------

# The root Extension row's Short/Definition used to be the generic
# "Extension" / "An Extension" placeholders; replace with the real
# short text and definition for this specific extension.
$elements.Range("K2").Value = "Adjudication Date"
$elements.Range("L2").Value = "Date on which the payment status of the claim was adjudicated"
